$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flat array of the new/updated data for rows 2-10 (columns A-T), 20 values per row.
$data = @(
  "ECs","Efna1","Epha7","ECs",2,0.6666666666666666,21.78783233333333,65.363497,0.9024488799587679,0.9024488799587679,3,1,0.3173666666666666,0.9520999999999999,0.1863268828340317,0.1863268828340316,6.91473172152222,62.23258549369999,0.1681504867197805,0.1681504867197804,
  "ECs","Efna1","Epha7","FAPs",2,0.6666666666666666,21.78783233333333,65.363497,0.9024488799587679,0.9024488799587679,2,0.6666666666666666,0.2780386666666667,0.8341160000000001,0.1632373009158609,0.1632373009158609,6.05785985151689,54.520738663652,0.1473133193790111,0.147313319379011,
  "ECs","Efna1","Epha7","sCs",2,0.6666666666666666,21.78783233333333,65.363497,0.9024488799587679,0.9024488799587679,3,1,1.107873666666667,3.323621,0.6504358162501074,0.6504358162501074,24.13816569584855,217.243491262637,0.5869850738599763,0.5869850738599763,
  "FAPs","Efna1","Epha7","ECs",3,1,1.359006333333333,4.077019,0.05628984676448105,0.05628984676448104,3,1,0.3173666666666666,0.9520999999999999,0.1863268828340317,0.1863268828340316,0.4313033099888888,3.8817297899,0.01048831168283106,0.01048831168283105,
  "FAPs","Efna1","Epha7","FAPs",3,1,1.359006333333333,4.077019,0.05628984676448105,0.05628984676448104,2,0.6666666666666666,0.2780386666666667,0.8341160000000001,0.1632373009158609,0.1632373009158609,0.3778563089115556,3.400706780204,0.009188602654801294,0.009188602654801292,
  "FAPs","Efna1","Epha7","sCs",3,1,1.359006333333333,4.077019,0.05628984676448105,0.05628984676448104,3,1,1.107873666666667,3.323621,0.6504358162501074,0.6504358162501074,1.505607329533222,13.550465965799,0.0366129324268487,0.03661293242684869,
  "sCs","Efna1","Epha7","ECs",3,1,0.9961713333333334,2.988514,0.04126127327675106,0.04126127327675105,3,1,0.3173666666666666,0.9520999999999999,0.1863268828340317,0.1863268828340316,0.3161515754888889,2.8453641794,0.007688084431420157,0.007688084431420154,
  "sCs","Efna1","Epha7","FAPs",3,1,0.9961713333333334,2.988514,0.04126127327675106,0.04126127327675105,2,0.6666666666666666,0.2780386666666667,0.8341160000000001,0.1632373009158609,0.1632373009158609,0.2769741492915556,2.492767343624001,0.006735378882048585,0.006735378882048582,
  "sCs","Efna1","Epha7","sCs",3,1,0.9961713333333334,2.988514,0.04126127327675106,0.04126127327675105,3,1,1.107873666666667,3.323621,0.6504358162501074,0.6504358162501074,1.103631987688222,9.932687889194002,0.02683780996328232,0.02683780996328232
)

$width = 20
$startRow = 2
$rowCount = $data.Length / $width

for ($i = 0; $i -lt $rowCount; $i++) {
  $r = $startRow + $i
  for ($j = 0; $j -lt $width; $j++) {
    $col = $j + 1
    $idx = $i * $width + $j
    $ws.Cells.Item($r, $col).Value = $data[$idx]
  }
}
